# Update to use lake-only data: point the Lacustrine RawDataPath (B5) at the
# new "Lakes_NEAP_20240723.shp" file instead of the old
# "Lakes_NEAP_20240808_NoOverlapWithALUM.shp" file, and move the active
# selection to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("B5").Value = "\\fs1-cbr.nexus.csiro.au\{ev-neap}\work\extent\processing\NEAP_intermediate\Lakes_NEAP_20240723.shp"

$ws.Range("B5").Select()
